$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new log row (7): "Xong giao diện" / "Fix lỗi login Token, http headers" ---

# A7: next day's date, formatted like A6 (copy A6's number format so it
# lands on the workbook's existing short-date style instead of creating a
# brand-new custom format).
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("A7").Value = 43531

# B7: author name, plain Times New Roman text (same font as the rest of the log).
$ws.Range("B7").Value = "Xong giao diện"
$ws.Range("B7").Font.Name = "Times New Roman"
$ws.Range("B7").Font.Size = 13

# C7:D7: merged, centered note describing the work done.
$ws.Range("C7:D7").Merge()
$ws.Range("C7").Value = "Fix lỗi login Token, http headers"
$ws.Range("C7:D7").Font.Name = "Times New Roman"
$ws.Range("C7:D7").Font.Size = 13
$ws.Range("C7:D7").HorizontalAlignment = -4108   # xlCenter

$ws.Range("C7:D7").Select()
